# Finished ID reading process
# - Remove the now-unused "Assets" sheet (and its Asset/OrchestratorAssetFolder
#   config rows) now that asset handling is no longer used by the process.
# - Extend the "Constants" sheet with the new Document Understanding /
#   Form Extractor related configuration entries, and the new output file
#   name setting, used by the finished ID reading process.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# --- Remove the "Assets" worksheet entirely -------------------------------
$assets = $wb.Worksheets.Item("Assets")
[void]$assets.Delete()

# --- Update the "Constants" sheet with the new configuration rows ---------
$ws = $wb.Worksheets.Item("Constants")

# Row 11 used to hold "Ids_directory_path" / "IDs" - that entry moves down to
# row 12 to make room for the new "Output_excel_file" setting at row 11.
$ws.Range("A12").Value = "Ids_directory_path"
$ws.Range("B12").Value = "IDs"

$ws.Range("A13").Value = "Document_Type_Id"
$ws.Range("B13").Value = "Scanari.IDs.ID"

$ws.Range("B14").Value = "https://du.uipath.com/svc/formextractor"
$ws.Range("A14").Value = "Form_Extractor_Endpoint"

$ws.Range("B15").Value = "JdvGeW9ZOkz5K4KWc20YEh/6NnSoXnDk/aQFUKESnyk8KGQD0u/rirxyof5LXHIKxidNUEg7Hcs4qq0WtuBKhg=="
$ws.Range("A15").Value = "Form_Extractor_Api_Key"

$ws.Range("A16").Value = "Success_processed_dir"
$ws.Range("A17").Value = "Error_processed_dir"

$ws.Range("B16").Value = "IDs/Success"
$ws.Range("B17").Value = "IDs/Error"

$ws.Range("A11").Value = "Output_excel_file"
$ws.Range("B11").Value = "Output.xlsx"

# Match the author's final selection/active cell on the Constants sheet.
$ws.Activate()
[void]$ws.Range("A11").Select()

# Keep the Settings sheet active tab / selection as in the source workbook.
$settings = $wb.Worksheets.Item("Settings")
[void]$settings.Range("C4").Select()

$ws.Activate()

$wb.Save()
